$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Grim Memorial"
$ws.Range("D4").Value = 1

$ws.Range("E6").Select()
